$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 421.75
$ws.Range("I9").Value = 412.16666
$ws.Range("J9").Value = 450.5
$ws.Range("K9").Value = 412.16666
$ws.Range("L9").Value = 450.5
$ws.Range("M9").Value = -243.16666
$ws.Range("N9").Value = -788.5
$ws.Range("H40").Value = 96921.91
$ws.Range("I40").Value = 335030.34
$ws.Range("J40").Value = 3749.0435
$ws.Range("K40").Value = 335030.34
$ws.Range("L40").Value = 3749.0435
$ws.Range("M40").Value = -334855.34
$ws.Range("N40").Value = -4099.0435
$ws.Range("H62").Value = 7427.5
$ws.Range("I62").Value = 7095.2
$ws.Range("K62").Value = 7095.2
$ws.Range("M62").Value = -6471.2
$ws.Range("H65").Value = 7427.5
$ws.Range("I65").Value = 7095.2
$ws.Range("K65").Value = 35476
$ws.Range("M65").Value = -32356
$ws.Range("H74").Value = 5004.8
$ws.Range("I74").Value = 4757.3
$ws.Range("K74").Value = 4757.3
$ws.Range("M74").Value = -3821.3
$ws.Range("H77").Value = 5004.8
$ws.Range("I77").Value = 4757.3
$ws.Range("K77").Value = 23786.5
$ws.Range("M77").Value = -19106.5
$ws.Range("H98").Value = 1592.475
$ws.Range("I98").Value = 1629.3334
$ws.Range("J98").Value = 155
$ws.Range("K98").Value = 1629.3334
$ws.Range("L98").Value = 155
$ws.Range("M98").Value = -131.3334
$ws.Range("N98").Value = -3151
$ws.Range("H101").Value = 706.63635
$ws.Range("I101").Value = 826.6
$ws.Range("K101").Value = 2479.8
$ws.Range("M101").Value = -857.8000000000002
$ws.Range("H122").Value = 1592.475
$ws.Range("I122").Value = 1629.3334
$ws.Range("J122").Value = 155
$ws.Range("K122").Value = 4888.0002
$ws.Range("L122").Value = 465
$ws.Range("M122").Value = -2438.0002
$ws.Range("N122").Value = -5365
$ws.Range("H129").Value = 3161.2334
$ws.Range("I129").Value = 1794.2
$ws.Range("K129").Value = 5382.6
$ws.Range("M129").Value = -382.6000000000004
$ws.Range("H132").Value = 3614.4878
$ws.Range("I132").Value = 3641.7222
$ws.Range("K132").Value = 10925.1666
$ws.Range("M132").Value = -8395.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3008.102
$ws.Range("I61").Value = 2987.6487
$ws.Range("K61").Value = 2987.6487
$ws.Range("M61").Value = -2775.6487
$ws.Range("H74").Value = 2840.6191
$ws.Range("J74").Value = 3825.4
$ws.Range("L74").Value = 3825.4
$ws.Range("N74").Value = -5573.4
$ws.Range("H77").Value = 2840.6191
$ws.Range("J77").Value = 3825.4
$ws.Range("L77").Value = 19127
$ws.Range("N77").Value = -27863
$ws.Range("H110").Value = 1625.2941
$ws.Range("I110").Value = 1180
$ws.Range("J110").Value = 3072.5
$ws.Range("K110").Value = 1180
$ws.Range("L110").Value = 3072.5
$ws.Range("M110").Value = 865
$ws.Range("N110").Value = -7162.5
$ws.Range("H133").Value = 154222
$ws.Range("J133").Value = 164999.75
$ws.Range("L133").Value = 164999.75
$ws.Range("N133").Value = -170059.75
$ws.Range("H134").Value = 55000
$ws.Range("J134").Value = 55000
$ws.Range("L134").Value = 55000
$ws.Range("N134").Value = -65140
$ws.Range("H136").Value = 3008.102
$ws.Range("I136").Value = 2987.6487
$ws.Range("K136").Value = 8962.946100000001
$ws.Range("M136").Value = -6412.946100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 98857.5
$ws.Range("J13").Value = 98857.5
$ws.Range("L13").Value = 98857.5
$ws.Range("N13").Value = -99193.5
$ws.Range("H99").Value = 2350.8333
$ws.Range("I99").Value = 1368.3334
$ws.Range("K99").Value = 1368.3334
$ws.Range("M99").Value = 129.6666
$ws.Range("H105").Value = 3462.5
$ws.Range("I105").Value = 3500.36
$ws.Range("K105").Value = 3500.36
$ws.Range("M105").Value = -1753.36
$ws.Range("H107").Value = 542.5135
$ws.Range("I107").Value = 486.66666
$ws.Range("K107").Value = 486.66666
$ws.Range("M107").Value = 1433.33334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 207.93333
$ws.Range("I7").Value = 102.125
$ws.Range("J7").Value = 328.85715
$ws.Range("K7").Value = 102.125
$ws.Range("L7").Value = 328.85715
$ws.Range("M7").Value = 10.875
$ws.Range("N7").Value = -554.85715
$ws.Range("H31").Value = 4698.3
$ws.Range("I31").Value = 3409.4092
$ws.Range("K31").Value = 3409.4092
$ws.Range("M31").Value = -3114.4092
$ws.Range("H34").Value = 4698.3
$ws.Range("I34").Value = 3409.4092
$ws.Range("K34").Value = 3409.4092
$ws.Range("M34").Value = -3207.4092
$ws.Range("H99").Value = 2092.7
$ws.Range("I99").Value = 1926.625
$ws.Range("J99").Value = 2757
$ws.Range("K99").Value = 1926.625
$ws.Range("L99").Value = 2757
$ws.Range("M99").Value = -428.625
$ws.Range("N99").Value = -5753
$ws.Range("H124").Value = 46159.6
$ws.Range("J124").Value = 46621.777
$ws.Range("L124").Value = 46621.777
$ws.Range("N124").Value = -51531.777
$ws.Range("H126").Value = 2092.7
$ws.Range("I126").Value = 1926.625
$ws.Range("J126").Value = 2757
$ws.Range("K126").Value = 5779.875
$ws.Range("L126").Value = 8271
$ws.Range("M126").Value = -3309.875
$ws.Range("N126").Value = -13211
$ws.Range("H134").Value = 2439.7307
$ws.Range("I134").Value = 1911.7646
$ws.Range("K134").Value = 5735.293799999999
$ws.Range("M134").Value = -3200.293799999999
$ws.Range("H141").Value = 1143632.2
$ws.Range("J141").Value = 1143632.2
$ws.Range("L141").Value = 1143632.2
$ws.Range("N141").Value = -1153992.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 94999.5
$ws.Range("J37").Value = 94999.5
$ws.Range("L37").Value = 284998.5
$ws.Range("N37").Value = -285222.5
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2685
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -1908
$ws.Range("H103").Value = 665
$ws.Range("J103").Value = 516.6667
$ws.Range("L103").Value = 1550.0001
$ws.Range("N103").Value = -3308.0001
$ws.Range("H117").Value = 2269.2
$ws.Range("J117").Value = 3477.3333
$ws.Range("L117").Value = 10431.9999
$ws.Range("N117").Value = -17315.9999
$ws.Range("H129").Value = 1790.8948
$ws.Range("J129").Value = 2941.3
$ws.Range("L129").Value = 8823.900000000001
$ws.Range("N129").Value = -18823.9
$ws.Range("H131").Value = 4124
$ws.Range("J131").Value = 3536.4
$ws.Range("L131").Value = 10609.2
$ws.Range("N131").Value = -20689.2
$ws.Range("H140").Value = 10002867
$ws.Range("I140").Value = 11113863
$ws.Range("K140").Value = 33341589
$ws.Range("M140").Value = -33336409

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 17503498
$ws.Range("I35").Value = 70000000
$ws.Range("K35").Value = 70000000
$ws.Range("M35").Value = -69999702
$ws.Range("H126").Value = 1367.6666
$ws.Range("I126").Value = 1347.75
$ws.Range("J126").Value = 1407.5
$ws.Range("K126").Value = 4043.25
$ws.Range("L126").Value = 4222.5
$ws.Range("M126").Value = -1573.25
$ws.Range("N126").Value = -9162.5
$ws.Range("H132").Value = 3337.6191
$ws.Range("I132").Value = 2671.7334
$ws.Range("K132").Value = 8015.2002
$ws.Range("M132").Value = -5485.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 124089.016
$ws.Range("I132").Value = 305947.3
$ws.Range("K132").Value = 917841.8999999999
$ws.Range("M132").Value = -915311.8999999999
$ws.Range("H136").Value = 3677.2666
$ws.Range("I136").Value = 3336.3845
$ws.Range("K136").Value = 10009.1535
$ws.Range("M136").Value = -7459.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 26214.695
$ws.Range("I136").Value = 2144.8
$ws.Range("K136").Value = 6434.400000000001
$ws.Range("M136").Value = -3884.400000000001
